# Adapt column header formatting to respective input file names.
#
# - The header row (row 1) columns used the generic "_old"/"_new" suffixes;
#   they are renamed to the input-file-specific "_FV2210"/"_FV2304" suffixes
#   (the "diff" column in between is left untouched).
# - The used range (A1:U93) is turned into a native Excel Table ("Table1")
#   so the new header names are also reflected in the table's column
#   definitions.
# - The top header row is frozen for easier scrolling through the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row cells -----------------------------------
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range into a native Excel table ------------------
$tableRange = $ws.Range("A1:U93")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
